$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = 175858
$ws.Range("C4").Value = 165830
$ws.Range("C7").Value = 5.7
$ws.Range("C8").Value = 64.65000000000001
